# Mise a jour de l'application
# Adds a new match day (N3J8) result for every player (row 2-28) and,
# for the two players who also played the reserve-team game (R2J6),
# records that match too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# N3J8 (column EA = minutes played, EB = T/R/NR/HG status)
$ws.Range("EA2").Value  = 90
$ws.Range("EB2").Value  = "T"

$ws.Range("EB3").Value  = "HG"

$ws.Range("EB4").Value  = "NR"

$ws.Range("EB5").Value  = "HG"

$ws.Range("EB6").Value  = "HG"

$ws.Range("EA7").Value  = 90
$ws.Range("EB7").Value  = "T"

$ws.Range("EB8").Value  = "HG"

$ws.Range("EA9").Value  = 90
$ws.Range("EB9").Value  = "T"

$ws.Range("EB10").Value = "NR"

$ws.Range("EA11").Value = 13
$ws.Range("EB11").Value = "R"

$ws.Range("EB12").Value = "HG"

$ws.Range("EB13").Value = "HG"

$ws.Range("EB14").Value = "HG"

$ws.Range("EA15").Value = 77
$ws.Range("EB15").Value = "T"

$ws.Range("EA16").Value = 90
$ws.Range("EB16").Value = "T"

$ws.Range("EA17").Value = 2
$ws.Range("EB17").Value = "R"

$ws.Range("EB18").Value = "HG"

$ws.Range("EA19").Value = 77
$ws.Range("EB19").Value = "T"

$ws.Range("EA20").Value = 88
$ws.Range("EB20").Value = "T"

$ws.Range("EB21").Value = "HG"

$ws.Range("EA22").Value = 90
$ws.Range("EB22").Value = "T"

$ws.Range("EB23").Value = "HG"

$ws.Range("EA24").Value = 90
$ws.Range("EB24").Value = "T"

$ws.Range("EA25").Value = 13
$ws.Range("EB25").Value = "R"

$ws.Range("EB26").Value = "HG"

$ws.Range("EA27").Value = 90
$ws.Range("EB27").Value = "T"

$ws.Range("EA28").Value = 90
$ws.Range("EB28").Value = "T"

# R2J6 (column JO = minutes played, JP = T/R/NR/HG status) for the
# two reserve-team players who also featured in that match.
$ws.Range("JO3").Value  = 90
$ws.Range("JP3").Value  = "T"

$ws.Range("JO26").Value = 90
$ws.Range("JP26").Value = "T"

# Restore the view state (active cell) recorded in the saved workbook.
$ws.Range("JU23").Select()
